$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reorder the member account-statement detail rows (B16:F41) so entries are
# grouped by worker (Yamadis Camargo, Gina Camargo, Olga Florez, Marcela Merino)
# across their mora periods, replacing the previous period-grouped ordering.
# Row 21 (Yamadis / periodo 1812) keeps its original position/value.

$ws.Range("B16").Value = "CC"
$ws.Range("C16").Value = "22815519"
$ws.Range("D16").Value = "YAMADIS CAMARGO MARQUEZ"
$ws.Range("E16").Value = "1905"
$ws.Range("F16").Value = 20833
$ws.Range("B17").Value = "CC"
$ws.Range("C17").Value = "22815519"
$ws.Range("D17").Value = "YAMADIS CAMARGO MARQUEZ"
$ws.Range("E17").Value = "1904"
$ws.Range("F17").Value = 31249
$ws.Range("B18").Value = "CC"
$ws.Range("C18").Value = "22815519"
$ws.Range("D18").Value = "YAMADIS CAMARGO MARQUEZ"
$ws.Range("E18").Value = "1903"
$ws.Range("F18").Value = 31249
$ws.Range("B19").Value = "CC"
$ws.Range("C19").Value = "22815519"
$ws.Range("D19").Value = "YAMADIS CAMARGO MARQUEZ"
$ws.Range("E19").Value = "1902"
$ws.Range("F19").Value = 31249
$ws.Range("B20").Value = "CC"
$ws.Range("C20").Value = "22815519"
$ws.Range("D20").Value = "YAMADIS CAMARGO MARQUEZ"
$ws.Range("E20").Value = "1901"
$ws.Range("F20").Value = 31249
$ws.Range("B21").Value = "CC"
$ws.Range("C21").Value = "22815519"
$ws.Range("D21").Value = "YAMADIS CAMARGO MARQUEZ"
$ws.Range("E21").Value = "1812"
$ws.Range("F21").Value = 31249
$ws.Range("B22").Value = "CC"
$ws.Range("C22").Value = "1049566193"
$ws.Range("D22").Value = "GINA MARCELA CAMARGO MONROY"
$ws.Range("E22").Value = "1905"
$ws.Range("F22").Value = 20833
$ws.Range("B23").Value = "CC"
$ws.Range("C23").Value = "1049566193"
$ws.Range("D23").Value = "GINA MARCELA CAMARGO MONROY"
$ws.Range("E23").Value = "1904"
$ws.Range("F23").Value = 31249
$ws.Range("B24").Value = "CC"
$ws.Range("C24").Value = "1049566193"
$ws.Range("D24").Value = "GINA MARCELA CAMARGO MONROY"
$ws.Range("E24").Value = "1903"
$ws.Range("F24").Value = 31249
$ws.Range("B25").Value = "CC"
$ws.Range("C25").Value = "1049566193"
$ws.Range("D25").Value = "GINA MARCELA CAMARGO MONROY"
$ws.Range("E25").Value = "1902"
$ws.Range("F25").Value = 31249
$ws.Range("B26").Value = "CC"
$ws.Range("C26").Value = "1049566193"
$ws.Range("D26").Value = "GINA MARCELA CAMARGO MONROY"
$ws.Range("E26").Value = "1901"
$ws.Range("F26").Value = 31249
$ws.Range("B27").Value = "CC"
$ws.Range("C27").Value = "1049566193"
$ws.Range("D27").Value = "GINA MARCELA CAMARGO MONROY"
$ws.Range("E27").Value = "1812"
$ws.Range("F27").Value = 31249
$ws.Range("B28").Value = "CC"
$ws.Range("C28").Value = "42365405"
$ws.Range("D28").Value = "OLGA ISABEL FLOREZ MEZA"
$ws.Range("E28").Value = "1905"
$ws.Range("F28").Value = 20833
$ws.Range("B29").Value = "CC"
$ws.Range("C29").Value = "42365405"
$ws.Range("D29").Value = "OLGA ISABEL FLOREZ MEZA"
$ws.Range("E29").Value = "1904"
$ws.Range("F29").Value = 31249
$ws.Range("B30").Value = "CC"
$ws.Range("C30").Value = "42365405"
$ws.Range("D30").Value = "OLGA ISABEL FLOREZ MEZA"
$ws.Range("E30").Value = "1903"
$ws.Range("F30").Value = 31249
$ws.Range("B31").Value = "CC"
$ws.Range("C31").Value = "42365405"
$ws.Range("D31").Value = "OLGA ISABEL FLOREZ MEZA"
$ws.Range("E31").Value = "1902"
$ws.Range("F31").Value = 31249
$ws.Range("B32").Value = "CC"
$ws.Range("C32").Value = "42365405"
$ws.Range("D32").Value = "OLGA ISABEL FLOREZ MEZA"
$ws.Range("E32").Value = "1901"
$ws.Range("F32").Value = 31249
$ws.Range("B33").Value = "CC"
$ws.Range("C33").Value = "42365405"
$ws.Range("D33").Value = "OLGA ISABEL FLOREZ MEZA"
$ws.Range("E33").Value = "1812"
$ws.Range("F33").Value = 31249
$ws.Range("B34").Value = "CC"
$ws.Range("C34").Value = "1095811770"
$ws.Range("D34").Value = "MARCELA MERIÑO OSPINO"
$ws.Range("E34").Value = "1905"
$ws.Range("F34").Value = 20833
$ws.Range("B35").Value = "CC"
$ws.Range("C35").Value = "1095811770"
$ws.Range("D35").Value = "MARCELA MERIÑO OSPINO"
$ws.Range("E35").Value = "1904"
$ws.Range("F35").Value = 31249
$ws.Range("B36").Value = "CC"
$ws.Range("C36").Value = "1095811770"
$ws.Range("D36").Value = "MARCELA MERIÑO OSPINO"
$ws.Range("E36").Value = "1903"
$ws.Range("F36").Value = 31249
$ws.Range("B37").Value = "CC"
$ws.Range("C37").Value = "1095811770"
$ws.Range("D37").Value = "MARCELA MERIÑO OSPINO"
$ws.Range("E37").Value = "1902"
$ws.Range("F37").Value = 31249
$ws.Range("B38").Value = "CC"
$ws.Range("C38").Value = "1095811770"
$ws.Range("D38").Value = "MARCELA MERIÑO OSPINO"
$ws.Range("E38").Value = "1901"
$ws.Range("F38").Value = 31249
$ws.Range("B39").Value = "CC"
$ws.Range("C39").Value = "1095811770"
$ws.Range("D39").Value = "MARCELA MERIÑO OSPINO"
$ws.Range("E39").Value = "1812"
$ws.Range("F39").Value = 31249
$ws.Range("B40").Value = "CC"
$ws.Range("C40").Value = "1095811770"
$ws.Range("D40").Value = "MARCELA MERIÑO OSPINO"
$ws.Range("E40").Value = "1811"
$ws.Range("F40").Value = 31249
$ws.Range("B41").Value = "CC"
$ws.Range("C41").Value = "1095811770"
$ws.Range("D41").Value = "MARCELA MERIÑO OSPINO"
$ws.Range("E41").Value = "1810"
$ws.Range("F41").Value = 31249
